$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Date header
Replace-Text "2023-10-16 Monday" "2023-10-17 Tuesday"

# Multiplication problems (processed top-to-bottom, matching document order,
# so that the "11x79=" -> "45x75=" change happens after "45x75=" -> "18x24=")
Replace-Text "89×34=" "47×63="
Replace-Text "96×78=" "38×65="
Replace-Text "42×40=" "98×27="
Replace-Text "41×80=" "70×18="
Replace-Text "75×41=" "57×77="
Replace-Text "44×80=" "24×42="
Replace-Text "42×12=" "19×54="
Replace-Text "12×38=" "82×35="
Replace-Text "82×69=" "24×74="
Replace-Text "76×46=" "59×58="
Replace-Text "84×51=" "30×47="
Replace-Text "60×49=" "71×39="
Replace-Text "71×74=" "38×91="
Replace-Text "28×66=" "93×99="
Replace-Text "32×29=" "50×99="
Replace-Text "17×49=" "64×55="
Replace-Text "45×75=" "18×24="
Replace-Text "38×45=" "20×23="
Replace-Text "63×42=" "22×74="
Replace-Text "29×81=" "16×45="
Replace-Text "53×99=" "52×39="
Replace-Text "11×79=" "45×75="
Replace-Text "21×28=" "71×25="
Replace-Text "59×44=" "75×53="
Replace-Text "60×85=" "81×95="
